$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 26 and 27 had their data entries (columns B through AC) swapped
# between each other, while the id values in column A (24 and 25) stay
# in place.
$range26 = $ws.Range("B26:AC26")
$range27 = $ws.Range("B27:AC27")

$values26 = $range26.Value2
$values27 = $range27.Value2

$range26.Value2 = $values27
$range27.Value2 = $values26
